$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("sku", "name", "quantity", "cost_per", "total_cost")

for ($r = 2; $r -le 11; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $headers[$c - 1]
    }
}
